# Adds three new "Test Data" country sheets (Russia, Finland, Hungary) at the
# end of the workbook, following the same template as the existing
# "Greece" market sheet, and moves tab-selection to the newly added,
# last sheet ("Hungary") - matching the author's commit:
#   "Added Test Data For Hungary/Russia/Finland Market"

$wb = $excel.ActiveWorkbook

# The existing "Greece" sheet is the closest structural template for the
# three new market sheets (same 18-row layout, same shared-string order).
$template = $wb.Worksheets.Item("Greece")

function Add-MarketSheet {
    param($Name, $Ticket, $Market)

    $last = $wb.Worksheets.Item($wb.Worksheets.Count)
    $template.Copy($null, $last)
    $newSheet = $wb.Worksheets.Item($last.Index + 1)
    $newSheet.Name = $Name

    # B2 = ticket id (e.g. "NGC-2929/T2898"), B4 = "<Country> Market"
    $newSheet.Range("B2").Value = $Ticket
    $newSheet.Range("B4").Value = $Market

    # Rows 3-5 grow to a 2-line height once the sheet is filled in.
    $newSheet.Range("A3:A5").RowHeight = 28.8

    # Column D narrows to fit the (now wrapped) labels.
    $newSheet.Columns.Item(4).ColumnWidth = 8.43

    return $newSheet
}

$russia  = Add-MarketSheet "Russia"  "NGC-2929/T2898" "Russia Market"
$finland = Add-MarketSheet "Finland" "NGC-3130/T2941" "Finland Market"
$hungary = Add-MarketSheet "Hungary" "NGC-3104/T2990" "Hungary Market"

# Restore each new sheet's own selection before moving on to the next one.
$russia.Range("P9").Select()
$finland.Range("A1:D18").Select()
$hungary.Range("A1:D18").Select()

# Hungary (the new last sheet) becomes the active / selected tab.
$hungary.Select()
